# Re-parse data & regenerate the q06_tasks_used survey-answer sheet.
# Rows are keyed by donor_id (col A), sorted ascending, with refreshed
# category answers (col B) and a new export timestamp (col D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$donorIds = @(
    "0ce5dd49",
    "2c1001cb",
    "37cc37bf",
    "43faa0b9",
    "4abe3e88",
    "50164f59",
    "5cf70f79",
    "5da96769",
    "6ca3e2f6",
    "790a4fcb",
    "802cc63a",
    "85c3ea4d",
    "942dfafb",
    "a2d65af2",
    "a46f1771",
    "ad58f9da",
    "c7d9a301",
    "ce8732ff",
    "d6f1d567",
    "da9326c9",
    "e09ca7bf",
    "ef53a641"
)

$categories = @(
    "Writing & professional communication|Coding - programming help|Study revision - exam prep|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Language practice or translation|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Study revision - exam prep|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Language practice or translation|Other",
    "Brainstorming & personal ideas - fun|Coding - programming help|Other",
    "Writing & professional communication|Coding - programming help|Language practice or translation|Other",
    "Other",
    "Writing & professional communication|Coding - programming help|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Other",
    "Writing & professional communication|Coding - programming help|Other",
    "Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Other",
    "Brainstorming & personal ideas - fun|Coding - programming help|Other",
    "Writing & professional communication|Coding - programming help|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Language practice or translation|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Language practice or translation|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Other",
    "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Other",
    "Brainstorming & personal ideas - fun|Coding - programming help|Other",
    "Writing & professional communication|Coding - programming help|Study revision - exam prep|Other"
)

$timestamp = 45854.65154966665
$surveyQuestion = "q06_tasks_used"

for ($i = 0; $i -lt $donorIds.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $donorIds[$i]
    $ws.Cells.Item($r, 2).Value = $categories[$i]
    $ws.Cells.Item($r, 3).Value = $surveyQuestion
    $ws.Cells.Item($r, 4).Value = $timestamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
